$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: "Realizar atendimento ao cliente" -> "Gerar orçamento"
#    (scope the Find to paragraph 3 only, so the identical phrase that
#     also appears inside the "Objetivo" paragraph stays untouched)
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Realizar atendimento ao cliente", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Gerar orçamento", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. "Evento" paragraph: "Cliente solicita atendimento" -> "Cliente solicita serviço"
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Cliente solicita atendimento", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Cliente solicita serviço", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Insert the new "Receber resposta do orçamento" block right after the
#    paragraph "Gera orçamento e guarda as informações em Orçamento."
#    (currently paragraph 10) and before "Aguarda resposta do orçamento
#    por parte do cliente." (currently paragraph 11).
#
#    We insert one blank paragraph plus five content paragraphs in a
#    single shot (using CR separators) and then go back to clear the
#    inherited list-numbering and apply the correct run formatting.
# ---------------------------------------------------------------------
$anchor = $d.Paragraphs(10).Range
$insPos = $anchor.End
$insRange = $d.Range($insPos, $insPos)
$insRange.InsertAfter("`rReceber resposta do orçamento`rEvento: Cliente confirma orçamento.`rObjetivo: Receber a resposta do cliente sobre o orçamento e salvar o orçamento.`rTrabalhadores Envolvidos: `rAtendente`r")

# The freshly inserted paragraphs are now #11 (blank) .. #16 (Atendente);
# paragraph #17 is the untouched "Aguarda resposta..." paragraph.
for ($i = 11; $i -le 16; $i++) {
    $np = $d.Paragraphs($i)
    $np.Range.ListFormat.RemoveNumbers()
    $np.Style = "Normal"
}

# --- Paragraph 12: "Receber resposta do orçamento" (bold) ---
$pB = $d.Paragraphs(12).Range
$pB.Font.Name = "Arial"
$pB.Font.NameBi = "Arial"
$pB.Font.Bold = $true

# --- Paragraph 13: "Evento: Cliente confirma orçamento." ---
$pC = $d.Paragraphs(13).Range
$cStart = $pC.Start
$rEvento = $d.Range($cStart, $cStart + 6)
$rEvento.Font.Name = "Arial"
$rEvento.Font.NameBi = "Arial"
$rEvento.Font.Bold = $true
$rRest = $d.Range($cStart + 6, $pC.End - 1)
$rRest.Font.Name = "Arial"
$rRest.Font.NameBi = "Arial"
$rRest.Font.Bold = $false
$rColon = $d.Range($cStart + 6, $cStart + 8)
$rColon.Font.Name = "Arial"
$rColon.Font.NameBi = "Arial"
$rColon.Font.Bold = $false
$rTexto = $d.Range($cStart + 8, $pC.End - 1)
$rTexto.Font.Name = "Arial"
$rTexto.Font.NameBi = "Arial"
$rTexto.Font.Bold = $false
$rTexto.Font.Color = 0

# --- Paragraph 14: "Objetivo: Receber a resposta do cliente sobre o orçamento e salvar o orçamento." ---
$pD = $d.Paragraphs(14).Range
$dStart = $pD.Start
$rObjetivo = $d.Range($dStart, $dStart + 8)
$rObjetivo.Font.Name = "Arial"
$rObjetivo.Font.NameBi = "Arial"
$rObjetivo.Font.Bold = $true
$rObjRest = $d.Range($dStart + 8, $pD.End - 1)
$rObjRest.Font.Name = "Arial"
$rObjRest.Font.NameBi = "Arial"
$rObjRest.Font.Bold = $false

# --- Paragraph 15: "Trabalhadores Envolvidos: " ---
$pE = $d.Paragraphs(15).Range
$eStart = $pE.Start
$rTrab = $d.Range($eStart, $eStart + 24)
$rTrab.Font.Name = "Arial"
$rTrab.Font.NameBi = "Arial"
$rTrab.Font.Bold = $true
$rTrabRest = $d.Range($eStart + 24, $pE.End - 1)
$rTrabRest.Font.Name = "Arial"
$rTrabRest.Font.NameBi = "Arial"
$rTrabRest.Font.Bold = $false

# --- Paragraph 16: "Atendente" ---
$pF = $d.Paragraphs(16).Range
$pF.Font.Name = "Arial"
$pF.Font.NameBi = "Arial"
$pF.Font.Bold = $false

Write-Output "block inserted"

# ---------------------------------------------------------------------
# 4. Last numbered item ("Guarda o orçamento em Ordem de serviço." then
#    "Informa ao cliente o número da ordem de serviço.") becomes a single
#    paragraph reading "armazenar o orçamento." and the trailing
#    "Informa ao cliente..." paragraph is removed entirely.
# ---------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
Write-Output $last.Range.Text
$secondLast = $d.Paragraphs($d.Paragraphs.Count - 1)
Write-Output $secondLast.Range.Text

# Delete the very last paragraph ("Informa ao cliente o número da ordem de serviço.")
$last.Range.Delete() | Out-Null

# Replace the text of what is now the last paragraph.
$tail = $d.Paragraphs($d.Paragraphs.Count).Range
$tail.Find.Execute("Guarda o orçamento em Ordem de serviço.", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "armazenar o orçamento.", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Relocate the hidden "_GoBack" bookmark from the "Evento" paragraph
#    to the very end of the document (end of the now-last paragraph).
# ---------------------------------------------------------------------
$endPara = $d.Paragraphs($d.Paragraphs.Count).Range
$bmPos = $endPara.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
